$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Add a new comment line before the "bioc.p" source-code paragraph
#    ("Paquetes de Bioconductor" section)
# -----------------------------------------------------------------
$biocParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Style.NameLocal -eq "Source Code" -and $candidate.Range.Text.StartsWith("bioc.p")) {
        $biocParagraph = $candidate
        break
    }
}

$insertPoint = $d.Range($biocParagraph.Range.Start, $biocParagraph.Range.Start)

$xmlComment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rStyle w:val="CommentTok"/></w:rPr><w:t xml:space="preserve"># nombres de los paquetes de Bioconductor que vamos a instalar</w:t></w:r><w:r><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xmlComment)

# -----------------------------------------------------------------
# 2) Remove the "PoiClaClu" / "glmpca" / "ggbeeswarm" lines from the
#    cran.p source-code block, keeping the comment about color
#    palettes and the following 'gridExtra' / 'colorspace' entries.
# -----------------------------------------------------------------
$r1 = $d.Content.Duplicate
$r1.Find.Execute("paletas de colores para los gráficos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$r2 = $d.Content.Duplicate
$r2.Find.Execute('"ggbeeswarm", ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$delRange = $d.Range($r1.End, $r2.End)
$delRange.Delete()
